$wb = $excel.ActiveWorkbook

# --- Update the email addresses on the per-person sheets ---
# Cell C4 on each of these sheets is an explicitly unlocked cell inside an
# otherwise protected sheet, so it can be edited directly without having to
# unprotect/reprotect the sheet. These values are also looked up (via
# shared formulas) by the "SOME" sheet, so updating them here cascades the
# new values there too.

$wsBarnaby = $wb.Worksheets.Item("Barnaby Barnes")
$wsBarnaby.Range("C4").Value = "b.barnes@learnpad.eu"

$wsSally = $wb.Worksheets.Item("Sally Shugar")
$wsSally.Range("C4").Value = "s.shugar@learnpad.eu"

# Recalculate so the formulas on the "SOME" sheet pick up the new values
$excel.Calculate()

# --- Update view/selection state ---
# Move the selection on the "Sally Shugar" sheet from C40 to C5
[void]$wsSally.Range("C5").Select()

# Make "Sally Shugar" the active tab (this also clears tabSelected on the
# previously active "SOME" sheet and updates the workbook's activeTab index)
$wsSally.Activate()
